$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1897
$ws.Range("F4").Value = 74
$ws.Range("F5").Value = 779
$ws.Range("F8").Value = 951
$ws.Range("F9").Value = 1639
$ws.Range("F10").Value = 1296
$ws.Range("F11").Value = 1579
$ws.Range("F13").Value = 1577
$ws.Range("F15").Value = 1713
$ws.Range("F17").Value = 1140
$ws.Range("F18").Value = 383
$ws.Range("F19").Value = 59
$ws.Range("F20").Value = 115
$ws.Range("F21").Value = 1888
$ws.Range("F22").Value = 262
$ws.Range("F23").Value = 828
$ws.Range("F26").Value = 1282
$ws.Range("F27").Value = 1085
$ws.Range("F28").Value = 85
$ws.Range("F30").Value = 1211
$ws.Range("F32").Value = 1190
$ws.Range("F33").Value = 1143
$ws.Range("F34").Value = 294
$ws.Range("F38").Value = 1714
$ws.Range("F39").Value = 16
$ws.Range("F41").Value = 2085
$ws.Range("F42").Value = 101
$ws.Range("F43").Value = 843
$ws.Range("F46").Value = 811
$ws.Range("F47").Value = 123

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 51
$ws.Range("F5").Value = 1507
$ws.Range("F7").Value = 4
$ws.Range("F8").Value = 2613
$ws.Range("F9").Value = 1228
$ws.Range("F10").Value = 420
$ws.Range("F11").Value = 733
$ws.Range("F12").Value = 267
$ws.Range("F13").Value = 43
$ws.Range("F14").Value = 80
$ws.Range("F21").Value = 327
$ws.Range("F25").Value = 19
$ws.Range("F29").Value = 258
$ws.Range("F31").Value = 233
$ws.Range("F34").Value = 61
$ws.Range("F35").Value = 31
$ws.Range("F42").Value = 67
$ws.Range("F45").Value = 69

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 272
$ws.Range("F5").Value = 2935
$ws.Range("F6").Value = 4712
$ws.Range("F7").Value = 154
$ws.Range("F9").Value = 602
$ws.Range("F10").Value = 791
$ws.Range("F11").Value = 486
$ws.Range("F12").Value = 428
$ws.Range("F13").Value = 1168
$ws.Range("F15").Value = 793

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1897
$ws.Range("F3").Value = 272
$ws.Range("F5").Value = 791
$ws.Range("F6").Value = 74
$ws.Range("F7").Value = 428
$ws.Range("F8").Value = 428
$ws.Range("F9").Value = 1168
$ws.Range("F11").Value = 951
$ws.Range("F12").Value = 1639
$ws.Range("F13").Value = 1296
$ws.Range("F14").Value = 1579
$ws.Range("F16").Value = 1577
$ws.Range("F17").Value = 267
$ws.Range("F18").Value = 80
$ws.Range("F19").Value = 1713
$ws.Range("F20").Value = 1140
$ws.Range("F22").Value = 793
$ws.Range("F23").Value = 793
$ws.Range("F24").Value = 1888
$ws.Range("F25").Value = 262
$ws.Range("F26").Value = 828
$ws.Range("F29").Value = 1283
$ws.Range("F30").Value = 327
$ws.Range("F31").Value = 1085
$ws.Range("F32").Value = 85
$ws.Range("F33").Value = 1211
$ws.Range("F35").Value = 1190
$ws.Range("F38").Value = 1143
$ws.Range("F39").Value = 294
$ws.Range("F43").Value = 1714
$ws.Range("F44").Value = 16
$ws.Range("F46").Value = 2085
$ws.Range("F47").Value = 101
$ws.Range("F48").Value = 843
$ws.Range("F50").Value = 811
$ws.Range("F51").Value = 123
$ws.Range("F53").Value = 69
